$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "26.507.92"
$ws.Range("E2").Value = "  +0.25%  "

Set-TextValue "D3" "1.837.18"
$ws.Range("E3").Value = "  -0.18%  "

Set-TextValue "D4" "1.002"
$ws.Range("E4").Value = "  +0.18%  "

Set-TextValue "D5" "258.11"
$ws.Range("E5").Value = "  -1.32%  "

$ws.Range("E6").Value = "  +0.12%  "

Set-TextValue "D7" "0.5228"
$ws.Range("E7").Value = "  +0.72%  "

Set-TextValue "D8" "0.3158"
$ws.Range("E8").Value = "  -3.45%  "

Set-TextValue "D9" "0.06774"
$ws.Range("E9").Value = "  -0.16%  "

Set-TextValue "D10" "18.65"
$ws.Range("E10").Value = "  -0.04%  "

Set-TextValue "D11" "0.7759"
$ws.Range("E11").Value = "  +0.84%  "

Set-TextValue "D12" "0.07756"
$ws.Range("E12").Value = "  +0.57%  "

Set-TextValue "D13" "1.823.51"
$ws.Range("E13").Value = "  -0.87%  "

Set-TextValue "D14" "87.64"
$ws.Range("E14").Value = "  -0.75%  "

Set-TextValue "D15" "4.996"
$ws.Range("E15").Value = "  -0.73%  "

$ws.Range("E16").Value = "  +0.20%  "

Set-TextValue "D17" "13.82"
$ws.Range("E17").Value = "  -0.73%  "

$ws.Range("E18").Value = "  +0.13%  "

Set-TextValue "D19" "0.000007909"
$ws.Range("E19").Value = "  -0.75%  "

Set-TextValue "D20" "26.542.32"
$ws.Range("E20").Value = "  +0.35%  "

Set-TextValue "D21" "2.069.09"
$ws.Range("E21").Value = "  +0.01%  "

Set-TextValue "D22" "4.592"
$ws.Range("E22").Value = "  +0.35%  "

Set-TextValue "D23" "5.950"
$ws.Range("E23").Value = "  -0.19%  "

Set-TextValue "D24" "9.289"
$ws.Range("E24").Value = "  -2.03%  "

Set-TextValue "D25" "142.57"
$ws.Range("E25").Value = "  -1.40%  "

Set-TextValue "D26" "2.210"
$ws.Range("E26").Value = "  -0.47%  "

Set-TextValue "D27" "1.671"
$ws.Range("E27").Value = "  +1.74%  "

Set-TextValue "D28" "16.86"
$ws.Range("E28").Value = "  -0.83%  "

Set-TextValue "D29" "111.78"
$ws.Range("E29").Value = "  +0.37%  "

Set-TextValue "D30" "4.159"
$ws.Range("E30").Value = "  -0.74%  "

Set-TextValue "D31" "0.08712"
$ws.Range("E31").Value = "  -0.22%  "

Set-TextValue "D32" "4.053"
$ws.Range("E32").Value = "  -1.92%  "

Set-TextValue "D33" "0.04868"
$ws.Range("E33").Value = "  +1.05%  "

$ws.Range("E34").Value = "  +0.30%  "

Set-TextValue "D35" "0.7178"
$ws.Range("E35").Value = "  +1.51%  "

Set-TextValue "D36" "2.862"
$ws.Range("E36").Value = "  +0.89%  "

Set-TextValue "D37" "3.085"
$ws.Range("E37").Value = "  +0.18%  "

Set-TextValue "D38" "2.225"
$ws.Range("E38").Value = "  +0.01%  "

Set-TextValue "D39" "0.01727"
$ws.Range("E39").Value = "  -1.81%  "

Set-TextValue "D40" "0.4803"
$ws.Range("E40").Value = "  -0.77%  "

Set-TextValue "D41" "0.8924"
$ws.Range("E41").Value = "  +0.21%  "

Set-TextValue "D42" "110.17"
$ws.Range("E42").Value = "  -1.06%  "

Set-TextValue "D43" "5.917"
$ws.Range("E43").Value = "  -2.72%  "

$ws.Range("E44").Value = "  +0.18%  "

Set-TextValue "D45" "7.609"
$ws.Range("E45").Value = "  -1.29%  "

Set-TextValue "D46" "0.4153"
$ws.Range("E46").Value = "  +0.25%  "

Set-TextValue "D47" "8.984"
$ws.Range("E47").Value = "  +0.23%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D48" "0.1228"
$ws.Range("E48").Value = "  +0.53%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D49" "0.05814"
$ws.Range("E49").Value = "  -0.88%  "

Set-TextValue "D50" "34.76"
$ws.Range("E50").Value = "  -0.74%  "

Set-TextValue "D51" "0.8906"
$ws.Range("E51").Value = "  +0.28%  "
